$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "x" markers in column J for rows 3,4,5,6,8 (new, previously absent)
$ws.Range("J3").Value = "x"
$ws.Range("J4").Value = "x"
$ws.Range("J5").Value = "x"
$ws.Range("J6").Value = "x"
$ws.Range("J8").Value = "x"

# Fill in row 24 with a new sensor entry: "gyrometer"
$ws.Range("A24").Value = "gyrometer"
$ws.Range("B24").Value = 1
$ws.Range("C24").Value = 5
$ws.Range("F24").Value = 2
$ws.Range("H24").Value = 1
$ws.Range("J24").Value = "x"

# Move the active cell selection from C21 to A10
$ws.Range("A10").Select()
